$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4, shifting rows 4-11 down to 5-12.
# This keeps the same style/layout as the existing rows (since Excel
# carries formatting from the row above on insert), matching the
# original row 3 style.
$ws.Rows("4").Insert()

# Fill in the (previously empty) row 3 with the new director.
$ws.Range("A3").Value = "Aarti Vellimedu"
$ws.Range("B3").Value = "Marketing Intern"

# Fill in the newly inserted row 4 with the second new director.
$ws.Range("A4").Value = "Elin Min"
$ws.Range("B4").Value = "Marketing Intern"

# Update the active selection to match the authored workbook (B3).
$ws.Range("B3").Select()
